$d = $word.ActiveDocument
$sec = $d.Sections(1)

# BTec_Logo-Orange images live in the headers: rename image2.jpg -> image1.jpg
for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $h = $sec.Headers.Item($i)
    if ($h.Exists) {
        $ishapes = $h.Range.InlineShapes
        for ($j = 1; $j -le $ishapes.Count; $j++) {
            $shp = $ishapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image1.jpg"
            }
        }
    }
}

# PearsonLogo images live in the footers: rename image1.png -> image2.png
for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $f = $sec.Footers.Item($i)
    if ($f.Exists) {
        $ishapes = $f.Range.InlineShapes
        for ($j = 1; $j -le $ishapes.Count; $j++) {
            $shp = $ishapes.Item($j)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image2.png"
            }
        }
    }
}
